# Update the table style used by the three data tables (slides 14-16)
# from the custom "Table_0" style to the built-in "No Style, No Grid"
# table style ({9B7DA487-C30B-4699-9EFA-BAC02272AE72}).

$p = $ppt.ActivePresentation

$newStyleId = "{9B7DA487-C30B-4699-9EFA-BAC02272AE72}"
$targetSlides = @(14, 15, 16)

foreach ($slideIndex in $targetSlides) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}
